# Insert a new column K (between existing J and L columns) on the active
# sheet, for every data row that doesn't already have it. The new cell
# gets the same formatting as the corresponding J cell in that row and
# contains a single space character " " (the same shared string already
# used by column L in these rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row ranges (1-based worksheet rows) that are missing column K.
# Rows 94-103 and 338-347 already contain a K cell and must be left alone.
$ranges = @(
    @(2, 93),
    @(104, 337),
    @(348, 391)
)

foreach ($range in $ranges) {
    $startRow = $range[0]
    $endRow = $range[1]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $jCell = $ws.Cells.Item($r, 10)
        $kCell = $ws.Cells.Item($r, 11)

        # Copy J's formatting (style) onto K, then set K's value.
        $jCell.Copy()
        $kCell.PasteSpecial(-4122)
        $kCell.Value = " "
    }
}

$excel.CutCopyMode = 0
